$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the three "Klasna et al." replication rows (rows 14-16),
# which shifts the rows below (Chauchard, Solaz, Arvate) up.
$ws.Range("A14:M16").EntireRow.Delete()

# The old selection pointed at D20, which no longer exists in the
# shrunk data range (now A1:M16); reset it to the top-left cell the
# way a freshly re-saved sheet would.
$ws.Range("A1").Select() | Out-Null

